$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,16,17,21,22,23,36,37,41,42,43,56,57,61,62,63,76,77,81,82,83,84,85,86,87,95,96,100,101,102,103,104,105,106,114,115,119,120,121,122,123,124,125,133,134,138,139,140,141,142,143,144,152,153,157,158,159,160,161,162,163,171,172,176,177,178,191,192,196,197,198,211,212,216,217,218,231,232)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

$rows2 = @(25,28,45,48,65,68,180,183,200,203,220,223)
foreach ($r in $rows2) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}
